{"js": "// Apply the four MN Policy Assumptions text edits using Word JS API (Office.js)\n// search-and-replace approach: find a unique substring, then replace its Range text.\n\nconst edits = [\n  {\n    find: \"natural gas CCS capacity\",\n    replace: \"natural gas non-peaker capacity\"\n  },\n  {\n    find: \"Assume 1,500 MW \",\n    replace: \"Assume 550 MW \"\n  },\n  {\n    find: \"storage added by 2034\",\n    replace: \"natural gas peaker capacity added between 2030 and 2034\"\n  },\n  {\n    find: \", equivalent to about 2% per year\",\n    replace: \", equivalent to about 2-2.5% per year\"\n  }\n];\n\nfor (const { find, replace } of edits) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: \"${find}\"`);\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Apply the four MN Policy Assumptions text edits using Word COM interop (Find/Replace).\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 0, $false, $replaceText, 2)\n}\n\nReplace-Text \"natural gas CCS capacity\" \"natural gas non-peaker capacity\"\nReplace-Text \"Assume 1,500 MW \" \"Assume 550 MW \"\nReplace-Text \"storage added by 2034\" \"natural gas peaker capacity added between 2030 and 2034\"\nReplace-Text \", equivalent to about 2% per year\" \", equivalent to about 2-2.5% per year\"\n"}
